$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-11 Wednesday" "2024-09-12 Thursday"
Replace-Text "85×57=4845" "85×75=6375"
Replace-Text "92×86=7912" "79×29=2291"
Replace-Text "36×98=3528" "60×74=4440"
Replace-Text "83×95=7885" "75×59=4425"
Replace-Text "49×73=3577" "71×80=5680"
Replace-Text "71×40=2840" "71×92=6532"
Replace-Text "50×85=4250" "88×79=6952"
Replace-Text "17×48=816" "12×73=876"
Replace-Text "33×65=2145" "15×47=705"
Replace-Text "26×94=2444" "63×49=3087"
Replace-Text "67×15=1005" "27×15=405"
Replace-Text "34×13=442" "37×80=2960"
Replace-Text "60×13=780" "21×19=399"
Replace-Text "68×89=6052" "79×79=6241"
Replace-Text "51×55=2805" "84×37=3108"
Replace-Text "69×69=4761" "57×68=3876"
Replace-Text "54×18=972" "81×65=5265"
Replace-Text "91×79=7189" "50×70=3500"
Replace-Text "67×89=5963" "45×59=2655"
Replace-Text "36×46=1656" "65×13=845"
Replace-Text "34×81=2754" "96×95=9120"
Replace-Text "27×37=999" "78×53=4134"
Replace-Text "59×56=3304" "84×75=6300"
Replace-Text "93×99=9207" "20×47=940"
Replace-Text "94×61=5734" "51×34=1734"

Write-Output "Done"
